# This workbook (a FHIR ValueSet "Metadata" export) is being refreshed to a
# newer generator run: version bump, status flip to draft, new date, the
# Contact property split into two rows (publisher contact + author contact),
# a new Jurisdiction property row inserted, and the Description row gets a
# value copied from the Title. The "Include from LOINC" sheet is unchanged
# in content (only its shared-string indices shift because of the above).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata"

# --- Make room for the extra "Jurisdiction" row ------------------------
# Clone the fully-formatted last row (15) down into a new row 16 first, so
# every row from 12..16 already carries the correct (bordered / wrapped)
# style before we shuffle the text values down into place.
$ws1.Range("A15:B15").Copy($ws1.Range("A16:B16"))

# Shift the old rows 12-15 content down into 13-16 (bottom-up so nothing
# gets clobbered before it's been read).
$ws1.Range("A16").Value = $ws1.Range("A15").Value2
$ws1.Range("B16").Value = $ws1.Range("B15").Value2

$ws1.Range("A15").Value = $ws1.Range("A14").Value2
$ws1.Range("B15").Value = $ws1.Range("B14").Value2

$ws1.Range("A14").Value = $ws1.Range("A13").Value2
$ws1.Range("B14").Value = $ws1.Range("B13").Value2

$ws1.Range("A13").Value = $ws1.Range("A12").Value2
$ws1.Range("B13").Value = $ws1.Range("B5").Value2   # Description <- Title text

# New row 12: Jurisdiction property (value left blank, as published)
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""

# --- Field updates -------------------------------------------------------
$ws1.Range("B3").Value  = "0.1.7"                                  # Version
$ws1.Range("B6").Value  = "draft"                                  # Status
$ws1.Range("B8").Value  = "2024-08-27T12:23:18-05:00"               # Date

# Contact: row 10 now carries the publisher's contact detail, row 11 the
# author's contact detail (previously both rows showed the same placeholder).
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"
